# Regenerate orders with updated distance/sizes.
# The experiment's Distance codes (D51 -> D55, D64 -> D69, D80 -> D86) and the
# "large" Size code (S30 -> S31) changed; every cell whose text contains one
# of those tokens (Condition, Filename_Left, Filename_Right, Distance, Size)
# needs the token swapped, everywhere it occurs, while leaving numeric /
# boolean columns (Trial, Duration_Seconds, Is_Repeat, Block, ConditionID)
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

$changed = 0

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2

        if ($val -is [string]) {
            $newVal = $val.Replace("D51", "D55").Replace("D64", "D69").Replace("D80", "D86").Replace("S30", "S31")
            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
                $changed = $changed + 1
            }
        }
    }
}

Write-Host "Cells changed: $changed"
